# Insert a new row for the "epitraxr" tool above the existing "epiworld" row,
# shifting all subsequent rows (and row-numbered references) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = "epitraxr: Manipulate Epitrax Data And Generate Reports"
$ws.Range("B21").Value = "R package for manipulating Epitrax data and generating reports. This tool is designed to simplify the process of working with Epitrax data, making it easier for public health professionals to analyze and report on disease surveillance data. It generates internal and public reports in CSV, Excel, or PDF formats."
$ws.Range("C21").Value = "Andrew Pulsipher"
$ws.Range("D21").Value = "a.pulsipher@utah.edu"
$ws.Range("E21").Value = "Yes"
$ws.Range("G21").Value = "Published"
$ws.Range("H21").Value = "MIT"
$ws.Range("I21").Value = "R"
$ws.Range("J21").Value = "Public Health Professionals"
$ws.Range("K21").Value = "Beginner"
$ws.Range("L21").Value = "Decision Support tool"
$ws.Range("M21").Value = "Epitrax data"
$ws.Range("N21").Value = "https://epiforesite.github.io/epitraxr/, https://github.com/EpiForeSITE/epitraxr"
$ws.Range("O21").Value = "https://github.com/EpiForeSITE/epitraxr"
